$wb = $excel.ActiveWorkbook

# --- Sheet "inputs": clear column D (header + values) ---
$wsInputs = $wb.Worksheets.Item("inputs")
$wsInputs.Range("D1:D7").ClearContents()

# --- Sheet "species_predvars": set dist_to_highways, population_density,
#     TotalInspections, days_fished (columns F, I, J, K) to FALSE for all
#     species rows (2-67) ---
$wsSpecies = $wb.Worksheets.Item("species_predvars")
$wsSpecies.Range("F2:F67").Value = $false
$wsSpecies.Range("I2:K67").Value = $false

# --- Selections / active sheet bookkeeping ---
$wsInputs.Range("D1:D7").Select() | Out-Null
$wsSpecies.Range("I5").Select() | Out-Null
$wsSpecies.Activate() | Out-Null
